$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 431.257796713125
$ws.Range("C2").Value = 511.961010133125
$ws.Range("D2").Value = 431.257796713125
$ws.Range("E2").Value = 511.961010133125

$ws.Range("B3").Value = 431.257796713125
$ws.Range("C3").Value = 514.48298555249994
$ws.Range("D3").Value = 431.257796713125
$ws.Range("E3").Value = 514.48298555249994

$ws.Range("B1:E3").Select()
